$p = $ppt.ActivePresentation
$nl = [char]10

# --- Slide 2 notes: drop the three paragraphs that were moved to slide 3 ---
$notes2 = $p.Slides.Item(2).NotesPage
$body2 = $notes2.Shapes.Item(2).TextFrame.TextRange
$body2.Text = "Vorschaubild des Museumskomplex" + $nl + ": Eröffnung vorraussichtlich Mai 2019"

# --- Slide 3 notes: tweak the title and append the paragraphs moved from slide 2 ---
$notes3 = $p.Slides.Item(3).NotesPage
$body3 = $notes3.Shapes.Item(2).TextFrame.TextRange
$body3.Text = "Luftbild mit Aufteilung" + $nl + `
    "" + $nl + `
    ": Museumskomplex baut sich zusammen aus HdbG, Bavariathek und dem österreischen Stadel( Depot) " + $nl + `
    ": Bavariathek fungiert als Bindeglied zwischen Stadel und HdbG " + $nl + `
    ": Verwaltung des gesamten Museumskomplexes in der Bavariathek" + $nl + `
    "ngen "

# --- Slide 7 notes: shorten first paragraph, drop the "Entwurf von Placeholder..." paragraph (moved to slide 9) ---
$notes6 = $p.Slides.Item(7).NotesPage
$body6 = $notes6.Shapes.Item(2).TextFrame.TextRange
$body6.Text = "Modellierung Grundriss bzw. äußere Hülle,  Detailierte Modellierung des Projektraums. (Später mehr)" + $nl + `
    "" + $nl + `
    "Eventueller Entwurf eines UI mit einer Game Engine zum Platzieren der erstellten Beispielassets (Stichwort: Planungstool). Wird nur gemacht, wenn noch Zeit übrig ist (unwahrscheinlich) " + $nl + `
    "" + $nl + `
    "Ansonsten ist das Projektziel eine fertige Blenderdatei für den Stakeholder (Hr. Wasweiß ich von der Bavariathek) "

# --- Slide 9 notes: rewrite the first two paragraphs, append the paragraph moved from slide 7 ---
$notes7 = $p.Slides.Item(9).NotesPage
$body7 = $notes7.Shapes.Item(2).TextFrame.TextRange
$body7.Text = "Der Projektraum soll detailierter dargestellt werden als die Restlichen Räume." + $nl + `
    "" + $nl + `
    "Hier sollen in Zukunft z.B. Projekte von Schülern vorgestellt, Plakate ausgehängt oder Präsentationen vorbereitet werden." + $nl + `
    "" + $nl + `
    "Entwurf von Placeholder Assets mit welchen die Bavariathek weiterarbeiten kann; Stühle, Tische, Podium mit Microfon?, Präsentationsmaterial (ggf. Platzhalter „Kunst)" + $nl + `
    ""
